$wb = $excel.ActiveWorkbook

# Original sheet order: "UserList", "Sorting", "UserList (2)"
# Target sheet order:   "UserList" (was "UserList (2)"), "UserList-1" (was "UserList"), "Sorting"

# Rename the original short "UserList" sheet out of the way first so the
# name "UserList" is free for the sheet moving into the first slot.
$wsOld = $wb.Worksheets.Item("UserList")
$wsOld.Name = "UserList-1"

# Move "UserList (2)" to be the very first sheet in the workbook.
$wsNew = $wb.Worksheets.Item("UserList (2)")
$wsNew.Move($wb.Worksheets.Item(1))

# Re-fetch (the moved sheet is still addressable by its old name) and rename it.
$wsNew = $wb.Worksheets.Item("UserList (2)")
$wsNew.Name = "UserList"

# Make the newly-relocated "UserList" sheet the active tab / selection.
$wsNew.Activate()
$wsNew.Range("F6").Select() | Out-Null
